$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 48
$ws1.Range("F5").Value = 1610
$ws1.Range("F6").Value = 3222
$ws1.Range("F7").Value = 752
$ws1.Range("F8").Value = 1969
$ws1.Range("F9").Value = 1890
$ws1.Range("F10").Value = 967
$ws1.Range("F11").Value = 338
$ws1.Range("F13").Value = 1567
$ws1.Range("F14").Value = 331
$ws1.Range("F16").Value = 56
$ws1.Range("F17").Value = 1378
$ws1.Range("F19").Value = 582
$ws1.Range("F20").Value = 278
$ws1.Range("F21").Value = 10363
$ws1.Range("F22").Value = 9555
$ws1.Range("F24").Value = 631
$ws1.Range("F25").Value = 1798
$ws1.Range("F26").Value = 127
$ws1.Range("F27").Value = 353

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 48
$ws4.Range("F7").Value = 1610
$ws4.Range("F8").Value = 3222
$ws4.Range("F9").Value = 752
$ws4.Range("F10").Value = 1969
$ws4.Range("F11").Value = 1890
$ws4.Range("F12").Value = 967
$ws4.Range("F13").Value = 338
$ws4.Range("F15").Value = 1567
$ws4.Range("F16").Value = 331
$ws4.Range("F19").Value = 56
$ws4.Range("F21").Value = 1378
$ws4.Range("F23").Value = 582
$ws4.Range("F24").Value = 278
$ws4.Range("F25").Value = 10363
$ws4.Range("F26").Value = 9555
$ws4.Range("F28").Value = 631
$ws4.Range("F29").Value = 1798
$ws4.Range("F32").Value = 127
$ws4.Range("F33").Value = 353
